$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 597.3182
$ws.Cells.Item(12, 9).Value = 571.3125
$ws.Cells.Item(12, 10).Value = 666.6667
$ws.Cells.Item(12, 11).Value = 571.3125
$ws.Cells.Item(12, 12).Value = 666.6667
$ws.Cells.Item(12, 13).Value = -401.3125
$ws.Cells.Item(12, 14).Value = -1006.6667
$ws.Cells.Item(33, 8).Value = 336
$ws.Cells.Item(33, 9).Value = 295.53845
$ws.Cells.Item(33, 10).Value = 467.5
$ws.Cells.Item(33, 11).Value = 295.53845
$ws.Cells.Item(33, 12).Value = 467.5
$ws.Cells.Item(33, 13).Value = -66.53845000000001
$ws.Cells.Item(33, 14).Value = -925.5
$ws.Cells.Item(97, 8).Value = 3046
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 3046
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 9138
$ws.Cells.Item(97, 14).Value = -10130
$ws.Cells.Item(98, 8).Value = 890.875
$ws.Cells.Item(98, 9).Value = 916.93335
$ws.Cells.Item(98, 10).Value = 500
$ws.Cells.Item(98, 11).Value = 916.93335
$ws.Cells.Item(98, 12).Value = 500
$ws.Cells.Item(98, 13).Value = 581.06665
$ws.Cells.Item(98, 14).Value = -3496
$ws.Cells.Item(101, 8).Value = 675.0625
$ws.Cells.Item(101, 9).Value = 579.6923
$ws.Cells.Item(101, 10).Value = 1088.3334
$ws.Cells.Item(101, 11).Value = 1739.0769
$ws.Cells.Item(101, 12).Value = 3265.0002
$ws.Cells.Item(101, 13).Value = -117.0769
$ws.Cells.Item(101, 14).Value = -6509.0002
$ws.Cells.Item(112, 8).Value = 66668650
$ws.Cells.Item(112, 9).Value = 166667600
$ws.Cells.Item(112, 10).Value = 2677.7778
$ws.Cells.Item(112, 11).Value = 500002800
$ws.Cells.Item(112, 12).Value = 8033.3334
$ws.Cells.Item(112, 13).Value = -500001692
$ws.Cells.Item(112, 14).Value = -10249.3334
$ws.Cells.Item(122, 8).Value = 890.875
$ws.Cells.Item(122, 9).Value = 916.93335
$ws.Cells.Item(122, 10).Value = 500
$ws.Cells.Item(122, 11).Value = 2750.80005
$ws.Cells.Item(122, 12).Value = 1500
$ws.Cells.Item(122, 13).Value = -300.8000499999998
$ws.Cells.Item(122, 14).Value = -6400
$ws.Cells.Item(132, 8).Value = 4291.595
$ws.Cells.Item(132, 9).Value = 4538.778
$ws.Cells.Item(132, 10).Value = 3846.6667
$ws.Cells.Item(132, 11).Value = 13616.334
$ws.Cells.Item(132, 12).Value = 11540.0001
$ws.Cells.Item(132, 13).Value = -11086.334
$ws.Cells.Item(132, 14).Value = -16600.0001
$ws.Cells.Item(137, 8).Value = 1925763.1
$ws.Cells.Item(137, 9).Value = 3573270.8
$ws.Cells.Item(137, 10).Value = 3670.8333
$ws.Cells.Item(137, 11).Value = 10719812.4
$ws.Cells.Item(137, 12).Value = 11012.4999
$ws.Cells.Item(137, 13).Value = -10717262.4
$ws.Cells.Item(137, 14).Value = -16112.4999
$ws.Cells.Item(138, 8).Value = 4880282
$ws.Cells.Item(138, 9).Value = 1970.4667
$ws.Cells.Item(138, 10).Value = 7694692.5
$ws.Cells.Item(138, 11).Value = 5911.4001
$ws.Cells.Item(138, 12).Value = 23084077.5
$ws.Cells.Item(138, 13).Value = -771.4000999999998
$ws.Cells.Item(138, 14).Value = -23094357.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2150.5833
$ws.Cells.Item(122, 9).Value = 1667.3334
$ws.Cells.Item(122, 10).Value = 2633.8333
$ws.Cells.Item(122, 11).Value = 5002.0002
$ws.Cells.Item(122, 12).Value = 7901.499899999999
$ws.Cells.Item(122, 13).Value = -2552.0002
$ws.Cells.Item(122, 14).Value = -12801.4999
$ws.Cells.Item(132, 8).Value = 36960.465
$ws.Cells.Item(132, 9).Value = 25683.404
$ws.Cells.Item(132, 10).Value = 63273.61
$ws.Cells.Item(132, 11).Value = 77050.212
$ws.Cells.Item(132, 12).Value = 189820.83
$ws.Cells.Item(132, 13).Value = -74520.212
$ws.Cells.Item(132, 14).Value = -194880.83

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 536.5
$ws.Cells.Item(22, 9).Value = 536.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 536.5
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -363.5
$ws.Cells.Item(105, 8).Value = 2857.7778
$ws.Cells.Item(105, 9).Value = 2840
$ws.Cells.Item(105, 10).Value = 3000
$ws.Cells.Item(105, 11).Value = 2840
$ws.Cells.Item(105, 12).Value = 3000
$ws.Cells.Item(105, 13).Value = -1093
$ws.Cells.Item(105, 14).Value = -6494

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(108, 8).Value = 40000
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 40000
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).Value = 40000
$ws.Cells.Item(108, 14).Value = -47680
$ws.Cells.Item(122, 8).Value = 2245.4783
$ws.Cells.Item(122, 9).Value = 1960.75
$ws.Cells.Item(122, 10).Value = 2556.0908
$ws.Cells.Item(122, 11).Value = 5882.25
$ws.Cells.Item(122, 12).Value = 7668.2724
$ws.Cells.Item(122, 13).Value = -3432.25
$ws.Cells.Item(122, 14).Value = -12568.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(64, 8).Value = 3294.75
$ws.Cells.Item(64, 9).Value = 1337.5555
$ws.Cells.Item(64, 10).Value = 4060.6086
$ws.Cells.Item(64, 11).Value = 4012.6665
$ws.Cells.Item(64, 12).Value = 12181.8258
$ws.Cells.Item(64, 13).Value = -3742.6665
$ws.Cells.Item(64, 14).Value = -12721.8258
$ws.Cells.Item(67, 8).Value = 3294.75
$ws.Cells.Item(67, 9).Value = 1337.5555
$ws.Cells.Item(67, 10).Value = 4060.6086
$ws.Cells.Item(67, 11).Value = 4012.6665
$ws.Cells.Item(67, 12).Value = 12181.8258
$ws.Cells.Item(67, 13).Value = -3076.6665
$ws.Cells.Item(67, 14).Value = -14053.8258
$ws.Cells.Item(113, 8).Value = 566.62
$ws.Cells.Item(113, 9).Value = 515.5417
$ws.Cells.Item(113, 10).Value = 613.7692
$ws.Cells.Item(113, 11).Value = 1546.6251
$ws.Cells.Item(113, 12).Value = 1841.3076
$ws.Cells.Item(113, 13).Value = 623.3749
$ws.Cells.Item(113, 14).Value = -6181.3076

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2314.8096
$ws.Cells.Item(97, 9).Value = 2317.8572
$ws.Cells.Item(97, 10).Value = 2308.7144
$ws.Cells.Item(97, 11).Value = 2317.8572
$ws.Cells.Item(97, 12).Value = 2308.7144
$ws.Cells.Item(97, 13).Value = -1821.8572
$ws.Cells.Item(97, 14).Value = -3300.7144
$ws.Cells.Item(102, 8).Value = 1717.7646
$ws.Cells.Item(102, 9).Value = 1694.8182
$ws.Cells.Item(102, 10).Value = 1759.8334
$ws.Cells.Item(102, 11).Value = 1694.8182
$ws.Cells.Item(102, 12).Value = 1759.8334
$ws.Cells.Item(102, 13).Value = -72.81819999999993
$ws.Cells.Item(102, 14).Value = -5003.8334
$ws.Cells.Item(113, 8).Value = 4150
$ws.Cells.Item(113, 9).Value = 800
$ws.Cells.Item(113, 10).Value = 7500
$ws.Cells.Item(113, 11).Value = 800
$ws.Cells.Item(113, 12).Value = 7500
$ws.Cells.Item(113, 13).Value = 1370
$ws.Cells.Item(113, 14).Value = -11840

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 792.93335
$ws.Cells.Item(16, 9).Value = 678
$ws.Cells.Item(16, 10).Value = 2402
$ws.Cells.Item(16, 11).Value = 678
$ws.Cells.Item(16, 12).Value = 2402
$ws.Cells.Item(16, 13).Value = -508
$ws.Cells.Item(16, 14).Value = -2742
$ws.Cells.Item(68, 8).Value = 1770
$ws.Cells.Item(68, 9).Value = 1700
$ws.Cells.Item(68, 10).Value = 1875
$ws.Cells.Item(68, 11).Value = 1700
$ws.Cells.Item(68, 12).Value = 1875
$ws.Cells.Item(68, 13).Value = -951
$ws.Cells.Item(68, 14).Value = -3373
$ws.Cells.Item(71, 8).Value = 1770
$ws.Cells.Item(71, 9).Value = 1700
$ws.Cells.Item(71, 10).Value = 1875
$ws.Cells.Item(71, 11).Value = 8500
$ws.Cells.Item(71, 12).Value = 9375
$ws.Cells.Item(71, 13).Value = -4756
$ws.Cells.Item(71, 14).Value = -16863
$ws.Cells.Item(93, 8).Value = 1134.3
$ws.Cells.Item(93, 9).Value = 1134.3
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 1134.3
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 113.7
$ws.Cells.Item(122, 8).Value = 3758.5757
$ws.Cells.Item(122, 9).Value = 4476.4614
$ws.Cells.Item(122, 10).Value = 3291.95
$ws.Cells.Item(122, 11).Value = 13429.3842
$ws.Cells.Item(122, 12).Value = 9875.849999999999
$ws.Cells.Item(122, 13).Value = -10979.3842
$ws.Cells.Item(122, 14).Value = -14775.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 6550
$ws.Cells.Item(62, 9).Value = 4680
$ws.Cells.Item(62, 10).Value = 9666.666999999999
$ws.Cells.Item(62, 11).Value = 4680
$ws.Cells.Item(62, 12).Value = 9666.666999999999
$ws.Cells.Item(62, 13).Value = -4056
$ws.Cells.Item(62, 14).Value = -10914.667
$ws.Cells.Item(65, 8).Value = 6550
$ws.Cells.Item(65, 9).Value = 4680
$ws.Cells.Item(65, 10).Value = 9666.666999999999
$ws.Cells.Item(65, 11).Value = 23400
$ws.Cells.Item(65, 12).Value = 48333.335
$ws.Cells.Item(65, 13).Value = -20280
$ws.Cells.Item(65, 14).Value = -54573.335
$ws.Cells.Item(96, 8).Value = 2300.5715
$ws.Cells.Item(96, 9).Value = 2750
$ws.Cells.Item(96, 10).Value = 2120.8
$ws.Cells.Item(96, 11).Value = 2750
$ws.Cells.Item(96, 12).Value = 2120.8
$ws.Cells.Item(96, 13).Value = -1377
$ws.Cells.Item(96, 14).Value = -4866.8
$ws.Cells.Item(126, 8).Value = 947.1667
$ws.Cells.Item(126, 9).Value = 936.8
$ws.Cells.Item(126, 10).Value = 999
$ws.Cells.Item(126, 11).Value = 2810.4
$ws.Cells.Item(126, 12).Value = 2997
$ws.Cells.Item(126, 13).Value = -340.3999999999996
$ws.Cells.Item(126, 14).Value = -7937
